$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells to update, in (address, new text value) pairs.
# Values are written as Text to preserve the original string cell type
# (these columns store numbers/percentages as text in the source sheet),
# then the number format is restored to General to match the original styling.
$updates = @{
    'D2' = '314.96'
    'E2' = '3.40%'
    'G2' = '10'
    'D3' = '35.62'
    'E3' = '-0.28%'
    'G3' = '10'
    'D4' = '5.144'
    'E4' = '1.15%'
    'G4' = '10'
    'D5' = '0.08106'
    'E5' = '2.90%'
    'G5' = '10'
    'D6' = '2.149'
    'E6' = '1.40%'
    'G6' = '10'
    'D7' = '8.016'
    'E7' = '1.29%'
    'G7' = '10'
    'D8' = '4.150'
    'E8' = '1.03%'
    'G8' = '10'
    'D9' = '0.9267'
    'E9' = '0.90%'
    'G9' = '10'
    'D10' = '0.1026'
    'E10' = '6.04%'
    'G10' = '10'
    'D11' = '0.1884'
    'E11' = '2.05%'
    'G11' = '10'
    'D12' = '0.09275'
    'E12' = '7.60%'
    'G12' = '10'
    'D13' = '0.03608'
    'E13' = '2.09%'
    'G13' = '10'
    'D14' = '0.09907'
    'E14' = '-0.32%'
    'G14' = '10'
    'D15' = '0.001441'
    'E15' = '-0.42%'
    'G15' = '10'
    'D16' = '0.005733'
    'E16' = '1.19%'
    'G16' = '10'
    'E17' = '0.44%'
    'G17' = '10'
    'D18' = '2.782'
    'E18' = '5.09%'
    'G18' = '10'
    'D19' = '0.3367'
    'E19' = '-1.94%'
    'G19' = '10'
    'E20' = '1.03%'
    'G20' = '10'
    'D21' = '5.140'
    'E21' = '-1.01%'
    'G21' = '10'
    'D22' = '0.2341'
    'E22' = '6.16%'
    'G22' = '10'
    'D23' = '0.04583'
    'E23' = '0.80%'
    'G23' = '10'
    'D24' = '0.001247'
    'E24' = '0.95%'
    'G24' = '10'
    'D25' = '0.004699'
    'E25' = '-7.07%'
    'G25' = '10'
    'D26' = '0.0001253'
    'E26' = '-21.81%'
    'G26' = '10'
    'D27' = '0.0004511'
    'E27' = '-5.05%'
    'G27' = '10'
    'G28' = '10'
    'G29' = '10'
    'G30' = '10'
    'G31' = '10'
    'G32' = '10'
    'G33' = '10'
    'G34' = '10'
    'G35' = '10'
    'G36' = '10'
    'G37' = '10'
    'G38' = '10'
    'D39' = '0.01956'
    'E39' = '6.39%'
    'G39' = '10'
    'D40' = '0.04875'
    'E40' = '3.38%'
    'G40' = '10'
    'D41' = '0.007803'
    'E41' = '3.37%'
    'G41' = '10'
    'D42' = '0.1390'
    'E42' = '-0.49%'
    'G42' = '10'
    'D43' = '0.007831'
    'E43' = '1.34%'
    'G43' = '10'
    'D44' = '0.002156'
    'E44' = '-2.26%'
    'G44' = '10'
    'D45' = '0.01161'
    'E45' = '5.39%'
    'G45' = '10'
    'D46' = '0.00006509'
    'E46' = '2.97%'
    'G46' = '10'
    'D47' = '0.00000000752'
    'E47' = '0.20%'
    'G47' = '10'
    'D48' = '35.84'
    'E48' = '-24.31%'
    'G48' = '10'
    'D49' = '0.001908'
    'E49' = '-4.64%'
    'G49' = '10'
    'D50' = '0.00002105'
    'E50' = '0.20%'
    'G50' = '10'
    'D51' = '0.0002005'
    'E51' = '0.20%'
    'G51' = '10'
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).NumberFormat = "@"
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Style = "Normal"
}

Write-Host "Updated $($updates.Count) cells."